$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / clearly-non-numeric assignments (kept as text naturally)
$ws.Range("D2").Value = "63.936.87"
$ws.Range("E2").Value = "  +5.55%  "
$ws.Range("D3").Value = "2.733.18"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  +5.89%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "2.747.52"
$ws.Range("E9").Value = "  +4.22%  "
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("E11").Value = "  +5.97%  "
$ws.Range("E12").Value = "  +4.02%  "
$ws.Range("E13").Value = "  +4.57%  "
$ws.Range("D14").Value = "3.228.24"
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "63.801.90"
$ws.Range("E16").Value = "  +5.38%  "
$ws.Range("E17").Value = "  +6.99%  "
$ws.Range("D18").Value = "2.746.49"
$ws.Range("E18").Value = "  +4.31%  "
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("E27").Value = "  +5.24%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "0.0₃0902"
$ws.Range("E29").Value = "  +12.99%  "
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("E31").Value = "  +6.37%  "
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  +14.67%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +4.41%  "
$ws.Range("E36").Value = "  +7.17%  "
$ws.Range("E37").Value = "  +9.72%  "
$ws.Range("E38").Value = "  +9.87%  "
$ws.Range("E39").Value = "  +13.55%  "
$ws.Range("E40").Value = "  +3.82%  "
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E43").Value = "  +8.30%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E46").Value = "  +5.00%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E47").Value = "  +6.01%  "
$ws.Range("E48").Value = "  +5.44%  "
$ws.Range("E49").Value = "  +4.38%  "
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("E51").Value = "  +0.11%  "

# Numeric-looking values that must stay as literal text (match source inlineStr),
# e.g. '1.00', '0.390', '21.75' would otherwise be auto-coerced to numbers by Excel
# and lose their exact textual formatting (trailing zeros, etc). Force text via a
# temporary "@" (Text) number format, then restore the default "Normal" style so no
# stray style/format is left behind on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.162"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.997"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.993"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "345.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.647"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0588"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"

Write-Output "done"
